$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stale "Output her er veldig midlertidig..." paragraph, along
#    with the blank paragraph that immediately followed it (the blank
#    paragraph before it, right after "Opptelling av typer..." is kept).
# ---------------------------------------------------------------------------
$marker = "Output her er veldig midlertidig"
$targetIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text.StartsWith($marker)) {
        $targetIndex = $idx
        break
    }
}

if ($targetIndex -gt 0) {
    $pOutput = $d.Paragraphs($targetIndex)
    $pNext = $d.Paragraphs($targetIndex + 1)
    $delStart = $pOutput.Range.Start
    $delEnd = $pNext.Range.End
    $d.Range($delStart, $delEnd).Delete()
}

# ---------------------------------------------------------------------------
# 2) Trim the "Input" paragraph: the sentence about Output being rewritten
#    ("Fordi Output er under omskriving, ...") is removed, joining the
#    remaining text back together.
# ---------------------------------------------------------------------------
$old = " for opptelling av registreringer per år. Fordi Output er under omskriving, kan dette også gjelde for denne delen av testen. Det kan også være nødvendig "
$new = " for opptelling av registreringer per år. Det kan også være nødvendig "

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2) | Out-Null
